$d = $word.ActiveDocument
$p2 = $d.Paragraphs.Item(2).Range
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="002A1F2A" w:rsidRDefault="002A1F2A" w:rsidP="002A1F2A"><w:pPr><w:tabs><w:tab w:val="left" w:pos="3119"/></w:tabs></w:pPr><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve">m:if </w:instrText></w:r><w:r><w:instrText xml:space="preserve">self.name </w:instrText></w:r><w:r><w:instrText>=</w:instrText></w:r><w:r><w:instrText xml:space="preserve">XXXX</w:instrText></w:r><w:r w:rsidR="00233314"><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r><w:r><w:t xml:space="preserve">    </w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>&lt;---</w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>Expression &quot;self.name =&quot; is invalid: missing expression</w:t></w:r></w:p>
'@
$p2.InsertXML($xml)
Write-Host "Field count:" $d.Fields.Count
Write-Host "Para2 text:" $d.Paragraphs.Item(2).Range.Text
